$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# so Excel does not auto-convert them (matching the original inline-string type).
$numericLookingCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D16", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D41", "D43", "D46", "D47", "D48", "D49")
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range('D2').Value = '46.297.91'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '2.613.85'
$ws.Range('E3').Value = '  +10.37%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '315.10'
$ws.Range('E5').Value = '  +5.60%  '
$ws.Range('D6').Value = '100.44'
$ws.Range('E6').Value = '  +3.65%  '
$ws.Range('E7').Value = '  +6.60%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.583'
$ws.Range('E9').Value = '  +16.63%  '
$ws.Range('D10').Value = '38.77'
$ws.Range('E10').Value = '  +14.55%  '
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  +7.65%  '
$ws.Range('D12').Value = '8.37'
$ws.Range('E12').Value = '  +19.43%  '
$ws.Range('D13').Value = '3.009.07'
$ws.Range('E13').Value = '  +10.26%  '
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '2.614.32'
$ws.Range('E15').Value = '  +10.39%  '
$ws.Range('D16').Value = '0.911'
$ws.Range('E16').Value = '  +12.03%  '
$ws.Range('D17').Value = '14.96'
$ws.Range('E17').Value = '  +9.59%  '
$ws.Range('D18').Value = '46.488.19'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('E19').Value = '  +8.72%  '
$ws.Range('D20').Value = '13.14'
$ws.Range('E20').Value = '  +3.94%  '
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  +12.14%  '
$ws.Range('D22').Value = '71.15'
$ws.Range('E22').Value = '  +6.92%  '
$ws.Range('D23').Value = '255.86'
$ws.Range('E23').Value = '  +6.11%  '
$ws.Range('D24').Value = '3.11'
$ws.Range('E24').Value = '  +13.42%  '
$ws.Range('D25').Value = '2.23'
$ws.Range('D26').Value = '28.29'
$ws.Range('E26').Value = '  +36.32%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '10.64'
$ws.Range('E28').Value = '  +10.87%  '
$ws.Range('D29').Value = '39.86'
$ws.Range('E29').Value = '  +3.93%  '
$ws.Range('E30').Value = '  +3.69%  '
$ws.Range('B31').Value = 'LidoDAOToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').Value = '3.75'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '6.15'
$ws.Range('E32').Value = '  +12.89%  '
$ws.Range('D33').Value = '2.31'
$ws.Range('E33').Value = '  +22.30%  '
$ws.Range('E34').Value = '  +5.66%  '
$ws.Range('D35').Value = '152.85'
$ws.Range('E35').Value = '  +4.31%  '
$ws.Range('D36').Value = '0.0835'
$ws.Range('E36').Value = '  +9.45%  '
$ws.Range('D37').Value = '0.118'
$ws.Range('E37').Value = '  +5.41%  '
$ws.Range('E38').Value = '  +6.01%  '
$ws.Range('D39').Value = '16.82'
$ws.Range('E39').Value = '  +11.07%  '
$ws.Range('D40').Value = '4.20'
$ws.Range('E40').Value = '  +10.42%  '
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').Value = '  +13.52%  '
$ws.Range('E42').Value = '  +11.47%  '
$ws.Range('D43').Value = '21.01'
$ws.Range('E43').Value = '  +49.45%  '
$ws.Range('D44').Value = '2.045.02'
$ws.Range('E44').Value = '  +5.48%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '91.29'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '9.29'
$ws.Range('E47').Value = '  +9.72%  '
$ws.Range('D48').Value = '110.25'
$ws.Range('E48').Value = '  +12.71%  '
$ws.Range('D49').Value = '1.79'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').Value = '2.867.08'
$ws.Range('E50').Value = '  +10.23%  '
$ws.Range('E51').Value = '  +10.03%  '

# Restore default style on the forced-text cells so only values changed
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).Style = "Normal"
}
